$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.902.97"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.872.90"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7396"
$ws.Range("E5").Value = "  -3.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.42"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3151"
$ws.Range("E8").Value = "  +0.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07218"
$ws.Range("E9").Value = "  +0.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.64"
$ws.Range("E10").Value = "  -3.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08314"
$ws.Range("E11").Value = "  -2.24%  "

$ws.Range("E12").Value = "  -1.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.379"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").Value = "1.858.02"
$ws.Range("E14").Value = "  -1.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.14"
$ws.Range("E15").Value = "  -1.66%  "

$ws.Range("D16").Value = "29.895.25"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "246.93"
$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("E19").Value = "  -1.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007829"
$ws.Range("E20").Value = "  +0.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").Value = "2.128.20"
$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.986"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1541"
$ws.Range("E25").Value = "  -4.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.283"
$ws.Range("E26").Value = "  -1.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.67"
$ws.Range("E27").Value = "  +2.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.64"
$ws.Range("E28").Value = "  -0.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.019"
$ws.Range("E29").Value = "  -0.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.495"
$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.571"
$ws.Range("E31").Value = "  +1.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.534"
$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.214"
$ws.Range("E33").Value = "  +2.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05319"

$ws.Range("E35").Value = "  -0.47%  "

$ws.Range("E36").Value = "  +0.82%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.698"
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  +0.70%  "

$ws.Range("E40").Value = "  -1.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4527"
$ws.Range("E41").Value = "  +1.41%  "

$ws.Range("D42").Value = "1.117.05"
$ws.Range("E42").Value = "  +1.26%  "

$ws.Range("E43").Value = "  +1.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.29"
$ws.Range("E44").Value = "  -0.81%  "

$ws.Range("E45").Value = "  +1.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.33"
$ws.Range("E46").Value = "  +1.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.860"
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.617"
$ws.Range("E49").Value = "  +0.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.490"
$ws.Range("E50").Value = "  -2.34%  "

$ws.Range("D51").Value = "2.028.72"
$ws.Range("E51").Value = "  +2.41%  "
